$wb = $excel.ActiveWorkbook

# --- Add new "Props" worksheet after "Datos" ---
$datos = $wb.Worksheets.Item("Datos")
$props = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $datos)
$props.Name = "Props"

$props.Range("A1").Value = "EA"
$props.Range("A2").Value = 1
$props.Range("B2").Value = 5875.2749999999996
$props.Range("A3").Value = 2
$props.Range("B3").Value = 200
$props.Range("A4").Value = 3
$props.Range("B4").Value = 300

# Copy the existing bordered cell style from the "Datos" sheet onto the new data
$datos.Range("A1").Copy()
$props.Range("A1").PasteSpecial(-4122)
$props.Range("A2:B4").PasteSpecial(-4122)

$props.Activate()
$excel.ActiveWindow.DisplayGridlines = $false
$props.Range("B5").Select() | Out-Null

# --- Update selection on "Datos" sheet ---
$datos.Activate()
$datos.Range("A3:B6").Select() | Out-Null

# --- Update selection on "Elementos" sheet, and leave it as the active tab ---
$elementos = $wb.Worksheets.Item("Elementos")
$elementos.Activate()
$elementos.Range("F17").Select() | Out-Null
